$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "living_rooms_2"
$ws.Range("E1").Value = "kitchens_2"
$ws.Range("F1").Value = "bedrooms_2"
